$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing-address line "2676 Greenrock Road, Milpitas CA 95035"
#    (the one that stands alone as the sender's return-address paragraph,
#    immediately followed by "Wayne Martinez") into two paragraphs:
#       "2676 Greenrock Road"
#       "Milpitas, CA 95035"
$rng = $d.Content
$rng.Find.Execute("2676 Greenrock Road, Milpitas CA 95035") | Out-Null
$rng.Text = "2676 Greenrock Road"
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("Milpitas, CA 95035")

# 3. Remove the empty "No Spacing" paragraph that immediately follows the
#    "... Board of Directors" signature line.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "Lees Orchard Association Board of Directors") {
        $next = $d.Paragraphs.Item($i + 1)
        $ntext = $next.Range.Text.TrimEnd([char]13, [char]7)
        if ($ntext -eq "" -and $next.Style.NameLocal -eq "No Spacing") {
            $next.Range.Delete()
        }
        break
    }
}
